$d = $word.ActiveDocument

# Original paragraph text is "Version 2." split across runs as:
#   "Versi" | "on" | " 2" | "." (with a spellcheck proofErr pair around
#   "Version" and a _GoBack bookmark between " 2" and ".")
# Target paragraph text is "Version 1." split as:
#   "Version" | " 1." (the "on" run is merged into the first run, and the
#   trailing "." run is removed, its text folded into the " 2" -> " 1." run)

# Step 1: remove the "on" run's text (characters 5-7 of "Versi" + "on").
$d.Range(5, 7).Text = ""

# Step 2: re-insert "on" right after "Versi" so it becomes part of the same
# run, producing a single run containing "Version".
$d.Range(5, 5).InsertAfter("on")

# Step 3: change "2" to "1." in place (stays inside the existing " 2" run,
# which becomes " 1.").
$d.Range(8, 9).Text = "1."

# Step 4: delete the now-redundant trailing "." run (characters 10-11).
$d.Range(10, 11).Text = ""
